$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.763.64'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '3.688.03'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('D4').Formula = '="2.41"'
$ws.Range('D4').Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('E4').Value = '  +27.54%  '
$ws.Range('D5').Formula = '="1.00"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').Formula = '="228.66"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  -3.41%  '
$ws.Range('D7').Formula = '="652.82"'
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Formula = '="0.439"'
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  +2.87%  '
$ws.Range('E9').Value = '  +8.24%  '
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').Value = '3.687.23'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Formula = '="47.98"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  +8.16%  '
$ws.Range('D13').Formula = '="0.209"'
$ws.Range('D13').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  +2.48%  '
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('D15').Formula = '="6.58"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  -2.86%  '
$ws.Range('D16').Value = '4.396.87'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').Value = '96.365.61'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Formula = '="8.89"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('D19').Value = '3.684.73'
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').Formula = '="19.11"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  +1.20%  '
$ws.Range('D21').Formula = '="12.93"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').Formula = '="0.543"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +7.70%  '
$ws.Range('D23').Formula = '="530.48"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  +2.44%  '
$ws.Range('D24').Formula = '="3.31"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  -1.83%  '
$ws.Range('D25').Formula = '="0.245"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +42.33%  '
$ws.Range('D26').Formula = '="119.40"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  +18.48%  '
$ws.Range('E27').Value = '  +1.74%  '
$ws.Range('D28').Formula = '="6.82"'
$ws.Range('D28').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('D29').Value = '3.886.79'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').Formula = '="12.89"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  -1.91%  '
$ws.Range('D31').Formula = '="13.39"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  +10.80%  '
$ws.Range('D33').Formula = '="1.00"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').Formula = '="0.187"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +1.52%  '
$ws.Range('D35').Formula = '="33.18"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  +3.03%  '
$ws.Range('D36').Formula = '="1.82"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  -2.05%  '
$ws.Range('D37').Formula = '="0.999"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').Value = '  +4.67%  '
$ws.Range('D39').Formula = '="612.62"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  -6.86%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').Formula = '="8.41"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  -5.01%  '
$ws.Range('D42').Formula = '="7.08"'
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  +3.26%  '
$ws.Range('E43').Value = '  +2.10%  '
$ws.Range('D44').Formula = '="0.0507"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +13.63%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Formula = '="0.488"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +12.46%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Formula = '="40.35"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  -2.49%  '
$ws.Range('D47').Formula = '="2.00"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  -4.06%  '
$ws.Range('D48').Formula = '="0.961"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('D49').Formula = '="8.98"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  +6.20%  '
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('D51').Formula = '="23.55"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  -0.06%  '
$excel.CutCopyMode = $false
